# Refresh the live cryptocurrency price/volume snapshot (cryptos.xlsx).
# Columns: A=rank(unchanged) B=Coin C=Link D=Price E=Volume(1h).
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as literal text (matching the sheet's inlineStr source)
# instead of silently coercing them to floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.747.20"
$ws.Range("E2").Value = "  -0.24%  "
# Row 3
$ws.Range("D3").Value = "2.025.33"
$ws.Range("E3").Value = "  -1.38%  "
# Row 4
$ws.Range("E4").Value = "  +0.15%  "
# Row 5
$ws.Range("D5").Value = "'226.72"
$ws.Range("E5").Value = "  -1.72%  "
# Row 6
$ws.Range("D6").Value = "'0.611"
$ws.Range("E6").Value = "  -0.79%  "
# Row 7
$ws.Range("D7").Value = "'59.49"
$ws.Range("E7").Value = "  +4.72%  "
# Row 8
$ws.Range("E8").Value = "  +0.06%  "
# Row 9
$ws.Range("D9").Value = "'0.383"
$ws.Range("E9").Value = "  -0.48%  "
# Row 10
$ws.Range("D10").Value = "'0.0806"
$ws.Range("E10").Value = "  +0.34%  "
# Row 11
$ws.Range("E11").Value = "  +0.26%  "
# Row 12
$ws.Range("D12").Value = "2.325.22"
$ws.Range("E12").Value = "  -1.17%  "
# Row 13
$ws.Range("D13").Value = "'14.50"
$ws.Range("E13").Value = "  -0.06%  "
# Row 14
$ws.Range("D14").Value = "'20.93"
$ws.Range("E14").Value = "  +2.08%  "
# Row 15
$ws.Range("D15").Value = "'0.749"
$ws.Range("E15").Value = "  +0.28%  "
# Row 16
$ws.Range("D16").Value = "'5.21"
$ws.Range("E16").Value = "  -0.82%  "
# Row 17
$ws.Range("D17").Value = "2.029.14"
$ws.Range("E17").Value = "  -1.13%  "
# Row 18
$ws.Range("D18").Value = "37.760.30"
$ws.Range("E18").Value = "  -0.02%  "
# Row 19
$ws.Range("D19").Value = "'6.04"
$ws.Range("E19").Value = "  -3.08%  "
# Row 20
$ws.Range("D20").Value = "'69.39"
$ws.Range("E20").Value = "  -0.26%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  -0.99%  "
# Row 22
$ws.Range("D22").Value = "'224.25"
$ws.Range("E22").Value = "  +0.02%  "
# Row 23
$ws.Range("E23").Value = "  +0.08%  "
# Row 24
$ws.Range("E24").Value = "  -1.75%  "
# Row 25
$ws.Range("D25").Value = "'2.22"
$ws.Range("E25").Value = "  -2.32%  "
# Row 26
$ws.Range("D26").Value = "'165.20"
$ws.Range("E26").Value = "  -0.23%  "
# Row 27
$ws.Range("D27").Value = "'9.16"
$ws.Range("E27").Value = "  -1.11%  "
# Row 28
$ws.Range("E28").Value = "  -3.59%  "
# Row 29
$ws.Range("D29").Value = "'18.82"
$ws.Range("E29").Value = "  -1.39%  "
# Row 30
$ws.Range("E30").Value = "  -5.33%  "
# Row 32
$ws.Range("E32").Value = "  -2.23%  "
# Row 33
$ws.Range("E33").Value = "  +0.64%  "
# Row 34
$ws.Range("D34").Value = "'0.0600"
$ws.Range("E34").Value = "  -2.15%  "
# Row 35
$ws.Range("D35").Value = "'4.47"
$ws.Range("E35").Value = "  -1.39%  "
# Row 36
$ws.Range("D36").Value = "'6.34"
$ws.Range("E36").Value = "  +6.59%  "
# Row 37
$ws.Range("E37").Value = "  -5.54%  "
# Row 38
$ws.Range("E38").Value = "  -1.68%  "
# Row 39
$ws.Range("E39").Value = "  +0.00%  "
# Row 40
$ws.Range("D40").Value = "1.535.88"
$ws.Range("E40").Value = "  +3.52%  "
# Row 41
$ws.Range("D41").Value = "'0.0216"
$ws.Range("E41").Value = "  -0.46%  "
# Row 42
$ws.Range("D42").Value = "'96.20"
$ws.Range("E42").Value = "  -2.10%  "
# Row 43
$ws.Range("D43").Value = "'16.56"
$ws.Range("E43").Value = "  -1.10%  "
# Row 44
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  -1.72%  "
# Row 45
$ws.Range("D45").Value = "'0.0917"
$ws.Range("E45").Value = "  -3.61%  "
# Row 46
$ws.Range("E46").Value = "  -1.93%  "
# Row 47
$ws.Range("D47").Value = "'2.96"
$ws.Range("E47").Value = "  +0.58%  "
# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -1.82%  "
# Row 49
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").Value = "'3.88"
$ws.Range("E49").Value = "  -5.85%  "
# Row 50
$ws.Range("D50").Value = "'7.06"
$ws.Range("E50").Value = "  -1.04%  "
# Row 51
$ws.Range("D51").Value = "2.215.18"
$ws.Range("E51").Value = "  -0.91%  "
